$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.072082666666667
$ws.Range("H2").Value = 12.216248
$ws.Range("I2").Value = 0.3979924983064649
$ws.Range("J2").Value = 0.3979924983064649
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 17.79587766666667
$ws.Range("N2").Value = 53.38763299999999
$ws.Range("O2").Value = 0.1793479316144739
$ws.Range("P2").Value = 0.179347931614474
$ws.Range("Q2").Value = 72.46628498455378
$ws.Range("R2").Value = 652.1965648609839
$ws.Range("S2").Value = 0.0713791313693415
$ws.Range("T2").Value = 0.07137913136934151

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.072082666666667
$ws.Range("H3").Value = 12.216248
$ws.Range("I3").Value = 0.3979924983064649
$ws.Range("J3").Value = 0.3979924983064649
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 72.39518466666667
$ws.Range("N3").Value = 217.185554
$ws.Range("O3").Value = 0.7296030503252249
$ws.Range("P3").Value = 0.7296030503252251
$ws.Range("Q3").Value = 294.7991766312658
$ws.Range("R3").Value = 2653.192589681392
$ws.Range("S3").Value = 0.2903765407709537
$ws.Range("T3").Value = 0.2903765407709538

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.072082666666667
$ws.Range("H4").Value = 12.216248
$ws.Range("I4").Value = 0.3979924983064649
$ws.Range("J4").Value = 0.3979924983064649
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.257112
$ws.Range("N4").Value = 0.771336
$ws.Range("O4").Value = 0.00259119028895291
$ws.Range("P4").Value = 0.00259119028895291
$ws.Range("Q4").Value = 1.046981318592
$ws.Range("R4").Value = 9.422831867328
$ws.Range("S4").Value = 0.001031274296687819
$ws.Range("T4").Value = 0.001031274296687819

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.072082666666667
$ws.Range("H5").Value = 12.216248
$ws.Range("I5").Value = 0.3979924983064649
$ws.Range("J5").Value = 0.3979924983064649
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.4631083333333333
$ws.Range("N5").Value = 1.389325
$ws.Range("O5").Value = 0.004667233797203165
$ws.Range("P5").Value = 0.004667233797203166
$ws.Range("Q5").Value = 1.885815416955555
$ws.Range("R5").Value = 16.9723387526
$ws.Range("S5").Value = 0.001857524039129256
$ws.Range("T5").Value = 0.001857524039129257

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4.072082666666667
$ws.Range("H6").Value = 12.216248
$ws.Range("I6").Value = 0.3979924983064649
$ws.Range("J6").Value = 0.3979924983064649
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.314158666666666
$ws.Range("N6").Value = 24.942476
$ws.Range("O6").Value = 0.08379059397414486
$ws.Range("P6").Value = 0.08379059397414489
$ws.Range("Q6").Value = 33.85594139444977
$ws.Range("R6").Value = 304.703472550048
$ws.Range("S6").Value = 0.03334802783035253
$ws.Range("T6").Value = 0.03334802783035255

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 6.159473666666667
$ws.Range("H7").Value = 18.478421
$ws.Range("I7").Value = 0.6020075016935351
$ws.Range("J7").Value = 0.6020075016935351
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 17.79587766666667
$ws.Range("N7").Value = 53.38763299999999
$ws.Range("O7").Value = 0.1793479316144739
$ws.Range("P7").Value = 0.179347931614474
$ws.Range("Q7").Value = 109.6132398630548
$ws.Range("R7").Value = 986.519158767493
$ws.Range("S7").Value = 0.1079688002451324
$ws.Range("T7").Value = 0.1079688002451325

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.159473666666667
$ws.Range("H8").Value = 18.478421
$ws.Range("I8").Value = 0.6020075016935351
$ws.Range("J8").Value = 0.6020075016935351
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 72.39518466666667
$ws.Range("N8").Value = 217.185554
$ws.Range("O8").Value = 0.7296030503252249
$ws.Range("P8").Value = 0.7296030503252251
$ws.Range("Q8").Value = 445.9162335478038
$ws.Range("R8").Value = 4013.246101930234
$ws.Range("S8").Value = 0.4392265095542712
$ws.Range("T8").Value = 0.4392265095542713

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.159473666666667
$ws.Range("H9").Value = 18.478421
$ws.Range("I9").Value = 0.6020075016935351
$ws.Range("J9").Value = 0.6020075016935351
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.257112
$ws.Range("N9").Value = 0.771336
$ws.Range("O9").Value = 0.00259119028895291
$ws.Range("P9").Value = 0.00259119028895291
$ws.Range("Q9").Value = 1.583674593384
$ws.Range("R9").Value = 14.253071340456
$ws.Range("S9").Value = 0.00155991599226509
$ws.Range("T9").Value = 0.001559915992265091

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.159473666666667
$ws.Range("H10").Value = 18.478421
$ws.Range("I10").Value = 0.6020075016935351
$ws.Range("J10").Value = 0.6020075016935351
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.4631083333333333
$ws.Range("N10").Value = 1.389325
$ws.Range("O10").Value = 0.004667233797203165
$ws.Range("P10").Value = 0.004667233797203166
$ws.Range("Q10").Value = 2.852503583980555
$ws.Range("R10").Value = 25.672532255825
$ws.Range("S10").Value = 0.002809709758073909
$ws.Range("T10").Value = 0.002809709758073909

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 6.159473666666667
$ws.Range("H11").Value = 18.478421
$ws.Range("I11").Value = 0.6020075016935351
$ws.Range("J11").Value = 0.6020075016935351
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 8.314158666666666
$ws.Range("N11").Value = 24.942476
$ws.Range("O11").Value = 0.08379059397414486
$ws.Range("P11").Value = 0.08379059397414489
$ws.Range("Q11").Value = 51.21084136782177
$ws.Range("R11").Value = 460.897572310396
$ws.Range("S11").Value = 0.05044256614379233
$ws.Range("T11").Value = 0.05044256614379234
